# Generate Report for Handoff
#
# A fresh handoff was generated, so every sheet's "in sync" status moves
# back to "Ready for handoff" and the handoff timestamps are refreshed.
# The (now shorter) status text also lets the status columns be narrowed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Reachable on the engine's column-width grid (1/6-character steps);
# 16.333333333333332 is the closest input to the authored 17.216 width.
$newStatusColWidth = 16.333333333333332

# --- Overview sheet ---
# E2/F2 hold the per-language handoff status, G2 the "Latest HO Xliff
# Generate Date" timestamp.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-24 07:00:58"

$wsOverview.Columns(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns(6).ColumnWidth = $newStatusColWidth

# --- zh-cn sheet ---
# C2 holds the handoff Status, H2 the Latest Handoff Datetime.
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-24 07:00:53"

$wsZhCn.Columns(3).ColumnWidth = $newStatusColWidth

# --- de-de sheet ---
# C2 holds the handoff Status, H2 the Latest Handoff Datetime (this is
# the timestamp mirrored by Overview!G2, since de-de was generated last).
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-24 07:00:58"

$wsDeDe.Columns(3).ColumnWidth = $newStatusColWidth
